$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.52877119971724
$ws.Range("D2").Value = 8.40347445535596
$ws.Range("E2").Value = 23.72370644997373
$ws.Range("F2").Value = 47.81549568951952
$ws.Range("G2").Value = 3.647214029169834
$ws.Range("L2").Value = 12.9215460309111
$ws.Range("N2").Value = 18.85471472147608
$ws.Range("B3").Value = 20.1997460090035
$ws.Range("D3").Value = 8.283939342072651
$ws.Range("E3").Value = 22.31678931846937
$ws.Range("F3").Value = 46.44910484407647
$ws.Range("G3").Value = 3.65626451365735
$ws.Range("L3").Value = 12.63759900783927
$ws.Range("N3").Value = 18.94955460255161
$ws.Range("B4").Value = 20.00154595407768
$ws.Range("D4").Value = 8.213834092094375
$ws.Range("E4").Value = 21.40617921535801
$ws.Range("F4").Value = 45.61402933852462
$ws.Range("G4").Value = 3.662078598073482
$ws.Range("L4").Value = 12.46454665655285
$ws.Range("N4").Value = 19.0100563128558
$ws.Range("B5").Value = 19.92183796870189
$ws.Range("D5").Value = 8.186118582628337
$ws.Range("E5").Value = 21.02341869416784
$ws.Range("F5").Value = 45.27525961152626
$ws.Range("G5").Value = 3.664513002440891
$ws.Range("L5").Value = 12.39445063803653
$ws.Range("N5").Value = 19.03528647337655
$ws.Range("B6").Value = 19.90866940148819
$ws.Range("D6").Value = 8.181568660057383
$ws.Range("E6").Value = 20.9591590110516
$ws.Range("F6").Value = 45.21911601606983
$ws.Range("G6").Value = 3.664921180349512
$ws.Range("L6").Value = 12.38283979161748
$ws.Range("N6").Value = 19.03951079528693
$ws.Range("B7").Value = 20.00046656249945
$ws.Range("D7").Value = 8.213456825392084
$ws.Range("E7").Value = 21.40106434016824
$ws.Range("F7").Value = 45.60945363444272
$ws.Range("G7").Value = 3.662111165008001
$ws.Range("L7").Value = 12.463599467467
$ws.Range("N7").Value = 19.01039424062209
$ws.Range("B8").Value = 20.41458971887971
$ws.Range("D8").Value = 8.36159212552421
$ws.Range("E8").Value = 23.24834449632139
$ws.Range("F8").Value = 47.34388076428927
$ws.Range("G8").Value = 3.650281576888771
$ws.Range("L8").Value = 12.82343200298176
$ws.Range("N8").Value = 18.88694718354571
$ws.Range("B9").Value = 21.25261361069282
$ws.Range("D9").Value = 8.67694675801741
$ws.Range("E9").Value = 26.49771386320795
$ws.Range("F9").Value = 50.75286705898464
$ws.Range("G9").Value = 3.629101248003533
$ws.Range("L9").Value = 13.53504975857958
$ws.Range("N9").Value = 18.66266768458383
$ws.Range("B10").Value = 21.87822169698809
$ws.Range("D10").Value = 8.922010676408116
$ws.Range("E10").Value = 28.65614122536401
$ws.Range("F10").Value = 53.23322908124866
$ws.Range("G10").Value = 3.614738636522294
$ws.Range("L10").Value = 14.05602518312536
$ws.Range("N10").Value = 18.50844322022803
$ws.Range("B11").Value = 22.16375285707459
$ws.Range("D11").Value = 9.03599751092295
$ws.Range("E11").Value = 29.58838773782352
$ws.Range("F11").Value = 54.35088500097046
$ws.Range("G11").Value = 3.608458073084156
$ws.Range("L11").Value = 14.29154473300158
$ws.Range("N11").Value = 18.44050856848247
$ws.Range("B12").Value = 22.27191196693516
$ws.Range("D12").Value = 9.079487225540987
$ws.Range("E12").Value = 29.93428267719133
$ws.Range("F12").Value = 54.77217118800178
$ws.Range("G12").Value = 3.606115616587553
$ws.Range("L12").Value = 14.38043407376711
$ws.Range("N12").Value = 18.41509785383337
$ws.Range("B13").Value = 22.24861798626407
$ws.Range("D13").Value = 9.070107021903983
$ws.Range("E13").Value = 29.86010470222064
$ws.Range("F13").Value = 54.68153235181037
$ws.Range("G13").Value = 3.606618520380004
$ws.Range("L13").Value = 14.36130465377289
$ws.Range("N13").Value = 18.4205565914113
$ws.Range("B14").Value = 22.17265103993087
$ws.Range("D14").Value = 9.039569114862335
$ws.Range("E14").Value = 29.61698745580489
$ws.Range("F14").Value = 54.3855856859205
$ws.Range("G14").Value = 3.608264642254281
$ws.Range("L14").Value = 14.29886412641347
$ws.Range("N14").Value = 18.43841173624822
$ws.Range("B15").Value = 22.12612068626625
$ws.Range("D15").Value = 9.020905101796711
$ws.Range("E15").Value = 29.46714350633387
$ws.Range("F15").Value = 54.20404485099917
$ws.Range("G15").Value = 3.609277593007307
$ws.Range("L15").Value = 14.26057642850587
$ws.Range("N15").Value = 18.44938935713652
$ws.Range("B16").Value = 21.85957180067711
$ws.Range("D16").Value = 8.914608804307871
$ws.Range("E16").Value = 28.5942186656289
$ws.Range("F16").Value = 53.15993685908701
$ws.Range("G16").Value = 3.615154132994953
$ws.Range("L16").Value = 14.04059638129249
$ws.Range("N16").Value = 18.51292721522842
$ws.Range("B17").Value = 21.69622070641336
$ws.Range("D17").Value = 8.850016428286441
$ws.Range("E17").Value = 28.04599977418735
$ws.Range("F17").Value = 52.51637223513656
$ws.Range("G17").Value = 3.618823638672938
$ws.Range("L17").Value = 13.90520599527231
$ws.Range("N17").Value = 18.55247144165765
$ws.Range("B18").Value = 21.60235966290239
$ws.Range("D18").Value = 8.813102588055839
$ws.Range("E18").Value = 27.72600577262856
$ws.Range("F18").Value = 52.14522845973172
$ws.Range("G18").Value = 3.62095809304993
$ws.Range("L18").Value = 13.82719894262459
$ws.Range("N18").Value = 18.57542576130939
$ws.Range("B19").Value = 21.57059921583409
$ws.Range("D19").Value = 8.80064613681968
$ws.Range("E19").Value = 27.61685854786493
$ws.Range("F19").Value = 52.01940996994744
$ws.Range("G19").Value = 3.621684895267852
$ws.Range("L19").Value = 13.8007667377874
$ws.Range("N19").Value = 18.5832338360242
$ws.Range("B20").Value = 21.71360069932867
$ws.Range("D20").Value = 8.856868020008534
$ws.Range("E20").Value = 28.10484220191281
$ws.Range("F20").Value = 52.58498556394762
$ws.Range("G20").Value = 3.618430548792337
$ws.Range("L20").Value = 13.91963304149242
$ws.Range("N20").Value = 18.54824024180638
$ws.Range("B21").Value = 22.1949642152415
$ws.Range("D21").Value = 9.048530286117455
$ws.Range("E21").Value = 29.68859016720652
$ws.Range("F21").Value = 54.47256829238326
$ws.Range("G21").Value = 3.607780167519791
$ws.Range("L21").Value = 14.31721312555448
$ws.Range("N21").Value = 18.4331587476237
$ws.Range("B22").Value = 22.50972116845561
$ws.Range("D22").Value = 9.175673132437328
$ws.Range("E22").Value = 30.68213305415438
$ws.Range("F22").Value = 55.69470706318057
$ws.Range("G22").Value = 3.60102827211476
$ws.Range("L22").Value = 14.57528940309803
$ws.Range("N22").Value = 18.35977806029605
$ws.Range("B23").Value = 22.34174721095565
$ws.Range("D23").Value = 9.10765404680212
$ws.Range("E23").Value = 30.15565447882279
$ws.Range("F23").Value = 55.04360731423833
$ws.Range("G23").Value = 3.604612964586807
$ws.Range("L23").Value = 14.43773729831933
$ws.Range("N23").Value = 18.39877678077319
$ws.Range("B24").Value = 21.70574303591789
$ws.Range("D24").Value = 8.853769723816844
$ws.Range("E24").Value = 28.0782545407373
$ws.Range("F24").Value = 52.55396902805878
$ws.Range("G24").Value = 3.61860818741996
$ws.Range("L24").Value = 13.91311109478761
$ws.Range("N24").Value = 18.55015248220247
$ws.Range("B25").Value = 21.02374180261979
$ws.Range("D25").Value = 8.589152012394875
$ws.Range("E25").Value = 25.65883071772554
$ws.Range("F25").Value = 49.83292485683926
$ws.Range("G25").Value = 3.634618312521695
$ws.Range("L25").Value = 13.3424795398714
$ws.Range("N25").Value = 18.72146612733381
